$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text format so numeric-looking price
# strings (e.g. "111.28") are written as text, matching the source data,
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "49.931.96"
$ws.Range("E2").Value = "  +3.80%  "

$ws.Range("D3").Value = "2.648.35"
$ws.Range("E3").Value = "  +5.85%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D5").Value = "111.28"
$ws.Range("E5").Value = "  +3.68%  "

$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "326.91"
$ws.Range("E6").Value = "  +2.13%  "

$ws.Range("E7").Value = "  +1.07%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "0.560"
$ws.Range("E9").Value = "  +3.45%  "

$ws.Range("D10").Value = "40.88"
$ws.Range("E10").Value = "  +2.96%  "

$ws.Range("D11").Value = "20.46"
$ws.Range("E11").Value = "  +1.66%  "

$ws.Range("D12").Value = "0.0822"
$ws.Range("E12").Value = "  +1.32%  "

$ws.Range("E13").Value = "  +0.59%  "

$ws.Range("D14").Value = "7.32"
$ws.Range("E14").Value = "  +3.24%  "

$ws.Range("D15").Value = "3.062.29"
$ws.Range("E15").Value = "  +5.84%  "

$ws.Range("D16").Value = "2.673.77"
$ws.Range("E16").Value = "  +6.81%  "

$ws.Range("D17").Value = "0.881"
$ws.Range("E17").Value = "  +5.67%  "

$ws.Range("D18").Value = "49.873.25"
$ws.Range("E18").Value = "  +3.96%  "

$ws.Range("D19").Value = "13.24"
$ws.Range("E19").Value = "  +2.24%  "

$ws.Range("E20").Value = "  +2.10%  "

$ws.Range("D21").Value = "2.95"
$ws.Range("E21").Value = "  +6.95%  "

$ws.Range("E22").Value = "  +2.21%  "

$ws.Range("D23").Value = "73.00"
$ws.Range("E23").Value = "  +2.03%  "

$ws.Range("D24").Value = "280.36"
$ws.Range("E24").Value = "  +1.31%  "

$ws.Range("E25").Value = "  +2.32%  "

$ws.Range("D26").Value = "27.04"
$ws.Range("E26").Value = "  +4.48%  "

$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("E28").Value = "  -1.93%  "

$ws.Range("D29").Value = "36.79"
$ws.Range("E29").Value = "  +4.80%  "

$ws.Range("D30").Value = "9.96"
$ws.Range("E30").Value = "  +2.48%  "

$ws.Range("E31").Value = "  +1.83%  "

$ws.Range("D32").Value = "49.85"
$ws.Range("E32").Value = "  +0.28%  "

$ws.Range("D33").Value = "19.62"
$ws.Range("E33").Value = "  +0.73%  "

$ws.Range("D34").Value = "5.45"
$ws.Range("E34").Value = "  +2.74%  "

$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("D36").Value = "0.0799"
$ws.Range("E36").Value = "  +2.05%  "

$ws.Range("D37").Value = "2.07"
$ws.Range("E37").Value = "  +6.58%  "

$ws.Range("D38").Value = "4.78"
$ws.Range("E38").Value = "  +3.04%  "

$ws.Range("E39").Value = "  +8.15%  "

$ws.Range("D40").Value = "127.40"
$ws.Range("E40").Value = "  +5.13%  "

$ws.Range("E41").Value = "  +1.68%  "

$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "22.49"
$ws.Range("E42").Value = "  +4.89%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "2.25"
$ws.Range("E43").Value = "  +1.74%  "

$ws.Range("E44").Value = "  +4.07%  "

$ws.Range("D45").Value = "3.39"
$ws.Range("E45").Value = "  +8.30%  "

$ws.Range("D46").Value = "2.068.06"
$ws.Range("E46").Value = "  +2.18%  "

$ws.Range("E47").Value = "  +13.96%  "

$ws.Range("E48").Value = "  +8.16%  "

$ws.Range("D49").Value = "9.09"
$ws.Range("E49").Value = "  +0.97%  "

$ws.Range("D50").Value = "5.41"
$ws.Range("E50").Value = "  +4.64%  "

$ws.Range("D51").Value = "82.01"
$ws.Range("E51").Value = "  +1.79%  "

# Restore default cell style on column D now that the text values are set,
# so the number format override does not linger on the cells.
$ws.Range("D2:D51").Style = "Normal"